$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the main worksheet from "Sheet1" to "centers"
$ws.Name = "centers"

# Restore all-caps column headers
$ws.Range("F1").Value = "BOROCODE"
$ws.Range("K1").Value = "ACCESSIBLE"
